$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G30").Value = "with_all_generated_comment"
$ws.Range("H30").Value = "eTour"
$ws.Range("I30").Value = "iTrust"
$ws.Range("J30").Value = "Albergate"

$ws.Range("G31").Value = "Correct"
$ws.Range("H31").Value = 30
$ws.Range("I31").Value = 136
$ws.Range("J31").Value = 32

$ws.Range("G32").Value = "Given"
$ws.Range("H32").Value = 68
$ws.Range("I32").Value = 501
$ws.Range("J32").Value = 122

$ws.Range("G33").Value = "Wanted"
$ws.Range("H33").Value = 308
$ws.Range("I33").Value = 418
$ws.Range("J33").Value = 54

$ws.Range("G34").Value = "Precision"
$ws.Range("H34").Value = "0.4412"
$ws.Range("I34").Value = "0.2715"
$ws.Range("J34").Value = "0.2623"

$ws.Range("G35").Value = "Recall"
$ws.Range("H35").Value = "0.0974"
$ws.Range("I35").Value = "0.3254"
$ws.Range("J35").Value = "0.5926"

$ws.Range("G36").Value = "F1"
$ws.Range("H36").Value = "0.1594"
$ws.Range("I36").Value = "0.2960"
$ws.Range("J36").Value = "0.3633"

$ws.Range("G34:J36").Font.Bold = $false

$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("G30:J36"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Tableau4"
$lo.TableStyle = "TableStyleLight12"
